# Adds a "preseason" column (K) with notes for phenology studies to the
# "decsens review" sheet.
#
# Commit message: "added preseaon information for phenology studies in
# full_review_guide spread sheet"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("decsens review")

# --- Column K: header + per-study preseason notes --------------------------
# Cells are written in the same order the source data was entered so that the
# generated shared-string table lines up with the authored workbook.

$ws.Range("K1").Value = "preseason"
$ws.Range("K1").Font.Bold = $true

$ws.Range("K2").Value = "no mention"

$ws.Range("K3").Value = "Nov1-April 30"

$ws.Range("K4").Value = "The optimal preseason for RP (relevant preseason) was defined as the period (at 15-day intervals) before the mean leaf unfolding date for which the Pearson correlation coefficient between leaf unfolding and air temperature was highest. The core algorithm for these models was ordinary least squares."
$ws.Range("K4").WrapText = $true
$ws.Range("K4").HorizontalAlignment = -4131
$ws.Range("K4").VerticalAlignment = -4160

$ws.Range("K5").Value = "Similarly to previous studies (Matsumoto et al. 2003; Dai et al. 2013), we estimated the optimum preseason (OP) that affects the phenophase most significantly for each species. First, we calculated the mean occurrence of LUD over the study period for each species, which was defined as the end date of OP (EP). We then calculated the start date of the OP (BP) by moving the date from EP-120 (120 days before EP) to EP-1 (the day before EP) in steps of 15 days. Pearson’s correlation analysis was performed to compare the series of data from each year and the mean temperature during each [BP, EP] period. The [BP, EP] period showing the correlation coefficient (R) with the highest absolute value was taken as the OP. (Dan's note 30-90 days was the OP for most)"
$ws.Range("K5").WrapText = $true

$ws.Range("K6").Value = " November to the mean date of leaf unfolding (MSOS)"
$ws.Range("K6").WrapText = $true

$ws.Range("K7").Value = "I think not applicable"

$ws.Range("K8").Value = "The phenology-based three-season division (no autumn) was determined from the partial correlation of the SOS to the climate of the individual month before the average SOS (Suppl. Material Fig. S1) and differs from a calendar-based four-season division."
$ws.Range("K8").WrapText = $true

$ws.Range("K11").Value = "to August of the year of tree-ring formation –"
$ws.Range("K11").WrapText = $true
$ws.Range("K11").VerticalAlignment = -4160

$ws.Range("K10").Value = "tarting from June of the previous year to August of the year of tree-ring formation"
$ws.Range("K10").WrapText = $true

$ws.Range("K9").Value = "specifically varied (one of the main points of this study)"
$ws.Range("K9").WrapText = $true
$ws.Range("K9").VerticalAlignment = -4160

$ws.Range("K25").Value = " The preseason at each station was defined as the period (with 5 day steps) before the mean FBD for which the Spearman's rank correlation coefficient between FBD and mean temperature was highest during 1963–2013 (H. Wang, Zhong, et al., 2017)."
$ws.Range("K25").WrapText = $true
$ws.Range("K25").VerticalAlignment = -4160

$ws.Range("K18").Value = "Here preseason temperature is defined as the mean temperature for the two months with the later month containing the 27-year average (1982–2008) of spring phenological date for each pixel"
$ws.Range("K18").WrapText = $true
$ws.Range("K18").VerticalAlignment = -4160

$ws.Range("K24").Value = "no mention"

# --- Row heights, adjusted to fit the new notes -----------------------------
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 180
$ws.Rows.Item(5).RowHeight = 409
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 150
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 75
$ws.Rows.Item(18).RowHeight = 165
$ws.Rows.Item(19).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(25).RowHeight = 150

# --- View state: leave the last-edited cell selected ------------------------
$ws.Range("K18").Select()

Write-Host "Added preseason column (K) with notes for phenology studies."
